$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha de Resultados")

$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2022-05-12"
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = "Comparação entre aumento da gasolina e patrimônio da família Bolsonaro usa dados imprecisos"
$ws.Range("C2").Value = "https://projetocomprova.com.br/publica%C3%A7%C3%B5es/comparacao-entre-aumento-da-gasolina-e-patrimonio-da-familia-bolsonaro-usa-dados-imprecisos/"

$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2022-05-12"
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Value = "Publicação que atribui a ex-tesoureiro do PT áudio contra igrejas é montagem"
$ws.Range("C3").Value = "https://projetocomprova.com.br/publica%C3%A7%C3%B5es/publicacao-que-atribui-a-ex-tesoureiro-do-pt-audio-contra-igrejas-e-montagem/"

$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "2022-05-10"
$ws.Range("A4").Style = "Normal"
$ws.Range("B4").Value = "É enganoso e está fora de contexto vídeo no TikTok em que Lula chama colaborador da Petrobras de corrupto"
$ws.Range("C4").Value = "https://projetocomprova.com.br/publica%C3%A7%C3%B5es/e-enganoso-e-esta-fora-de-contexto-video-no-tiktok-em-que-lula-chama-colaborador-da-petrobras-de-corrupto/"

$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "2022-05-09"
$ws.Range("A5").Style = "Normal"
$ws.Range("B5").Value = "Jovem que faz sátira sobre militantes de esquerda não é filha da deputada Maria do Rosário"
$ws.Range("C5").Value = "https://projetocomprova.com.br/publica%C3%A7%C3%B5es/jovem-que-faz-satira-sobre-militantes-de-esquerda-nao-e-filha-da-deputada-maria-do-rosario/"

$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "2022-05-06"
$ws.Range("A6").Style = "Normal"
$ws.Range("B6").Value = "Publicação tira de contexto declarações de Djavan para atacar a Lei Rouanet"
$ws.Range("C6").Value = "https://projetocomprova.com.br/publica%C3%A7%C3%B5es/publicacao-tira-de-contexto-declaracoes-de-djavan-para-atacar-a-lei-rouanet/"

$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "2022-05-05"
$ws.Range("A7").Style = "Normal"
$ws.Range("B7").Value = "Protesto de indígenas na Bahia era por melhoria na educação e não por verba para ato contra Bolsonaro"
$ws.Range("C7").Value = "https://projetocomprova.com.br/publica%C3%A7%C3%B5es/protesto-de-indigenas-na-bahia-era-por-melhoria-na-educacao-e-nao-por-verba-para-ato-contra-bolsonaro/"

$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "2022-05-04"
$ws.Range("A8").Style = "Normal"
$ws.Range("B8").Value = "Vídeos antigos são usados para enganar sobre adesão a atos pró-Bolsonaro em 1º de Maio"
$ws.Range("C8").Value = "https://projetocomprova.com.br/publica%C3%A7%C3%B5es/videos-antigos-sao-usados-para-enganar-sobre-adesao-a-atos-pro-bolsonaro-em-1o-de-maio/"

$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "2022-05-03"
$ws.Range("A9").Style = "Normal"
$ws.Range("B9").Value = "Antes de ser preso, João de Deus compareceu à posse de Rosa Weber no TSE"
$ws.Range("C9").Value = "https://projetocomprova.com.br/publica%C3%A7%C3%B5es/antes-de-ser-preso-joao-de-deus-compareceu-a-posse-de-rosa-weber-no-tse/"

$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "2022-05-03"
$ws.Range("A10").Style = "Normal"
$ws.Range("B10").Value = "É falso que Elon Musk tenha citado as motociatas de Bolsonaro em entrevista na Alemanha"
$ws.Range("C10").Value = "https://projetocomprova.com.br/publica%C3%A7%C3%B5es/e-falso-que-elon-musk-tenha-citado-as-motociatas-de-bolsonaro-em-entrevista-na-alemanha/"

$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = "2022-04-28"
$ws.Range("A11").Style = "Normal"
$ws.Range("B11").Value = "Post confunde dados e engana sobre conflitos no campo no governo Bolsonaro"
$ws.Range("C11").Value = "https://projetocomprova.com.br/publica%C3%A7%C3%B5es/post-confunde-dados-e-engana-sobre-conflitos-no-campo-no-governo-bolsonaro/"
